$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing cells
$ws.Range("G4").Value = 'Uploaded'
$ws.Range("G11").Value = 'Uploaded'

# Append new rows 12-23
# Row 12
$ws.Range("A12").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Slime Squishing\152BA.mp4"'
$ws.Range("B12").Value = 'Slime Squishing'
$ws.Range("C12").Value = '💘 ASMR Slime Squishing so Satisfying Relaxing #SlimeSquishing #slimeasmr #asmrslime #satisfying'
$ws.Range("D12").Value = '🍡💘 ASMR Slime Squishing so Satisfying Nutella 🍩
#SlimeSquishing #Slime #ASMR #Satisfying #Shorts #YouTubeShorts
#pipingbags  #slimeasmr #asmrslime 
slime asmr, asmr slime, satisfying
Thank you for watching the video. Please subscribe to Slime Squishing to watch our latest videos! Wish you have relaxing moments with Slime Squishing.'
$ws.Range("E12").Value = '22:00'
$ws.Range("F12").Value = '''11/9/2026'
$ws.Range("G12").Value = 'Uploaded'

# Row 13
$ws.Range("A13").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Slime Squishing\154BA.mp4"'
$ws.Range("B13").Value = 'Slime Squishing'
$ws.Range("C13").Value = '❤️🧡💛💚💙💜 ASMR Slime Squishing 💝 so Satisfying #SlimeSquishing #slimeasmr #asmrslime #satisfying'
$ws.Range("D13").Value = '❤️🧡💛💚💙💜 ASMR Slime Squishing 💝 so Satisfying
#SlimeSquishing #Slime #ASMR #Satisfying #Shorts #YouTubeShorts
#pipingbags  #slimeasmr #asmrslime 
slime asmr, asmr slime, satisfying
Thank you for watching the video. Please subscribe to Slime Squishing to watch our latest videos! Wish you have relaxing moments with Slime Squishing.'
$ws.Range("E13").Value = '2:00'
$ws.Range("F13").Value = '''12/9/2025'
$ws.Range("G13").Value = 'Uploaded'

# Row 14
$ws.Range("A14").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Slime Squishing\156BA.mp4"'
$ws.Range("B14").Value = 'Slime Squishing'
$ws.Range("C14").Value = '🍡💘 ASMR Slime Squishing so Satisfying Relaxing #SlimeSquishing #slimeasmr #asmrslime #satisfying'
$ws.Range("D14").Value = '🍡💘 ASMR Slime Squishing so Satisfying Nutella 🍩
#SlimeSquishing #Slime #ASMR #Satisfying #Shorts #YouTubeShorts
#pipingbags  #slimeasmr #asmrslime 
slime asmr, asmr slime, satisfying
Thank you for watching the video. Please subscribe to Slime Squishing to watch our latest videos! Wish you have relaxing moments with Slime Squishing.'
$ws.Range("E14").Value = '10:00'
$ws.Range("F14").Value = '''12/9/2025'
$ws.Range("G14").Value = 'Uploaded'

# Row 15
$ws.Range("A15").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Slime Squishing\158BA.mp4"'
$ws.Range("B15").Value = 'Slime Squishing'
$ws.Range("C15").Value = '💦 Topping slime 🧼 dried floam 💥 mixing ASMR  #SlimeSquishing #slimeasmr #asmrslime #satisfying'
$ws.Range("D15").Value = '💦 Topping slime 🧼 dried floam 💥 mixing ASMR
#SlimeSquishing #Slime #ASMR #Satisfying #Shorts #YouTubeShorts
#pipingbags  #slimeasmr #asmrslime 
slime asmr, asmr slime, satisfying
Thank you for watching the video. Please subscribe to Slime Squishing to watch our latest videos! Wish you have relaxing moments with Slime Squishing.'
$ws.Range("E15").Value = '22:00'
$ws.Range("F15").Value = '''12/9/2025'
$ws.Range("G15").Value = 'Uploaded'

# Row 16
$ws.Range("A16").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Slime Squishing\160BA.mp4"'
$ws.Range("B16").Value = 'Slime Squishing'
$ws.Range("C16").Value = '❤️🧡💛💚💙💜 ASMR Slime Squishing 💝 so Satisfying #SlimeSquishing #slimeasmr #asmrslime #satisfying'
$ws.Range("D16").Value = '❤️🧡💛💚💙💜 ASMR Slime Squishing 💝 so Satisfying
#SlimeSquishing #Slime #ASMR #Satisfying #Shorts #YouTubeShorts
#pipingbags  #slimeasmr #asmrslime 
slime asmr, asmr slime, satisfying
Thank you for watching the video. Please subscribe to Slime Squishing to watch our latest videos! Wish you have relaxing moments with Slime Squishing.'
$ws.Range("E16").Value = '2:00'
$ws.Range("F16").Value = '''13/9/2025'
$ws.Range("G16").Value = 'Uploaded'

# Row 17
$ws.Range("A17").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Slime Squishing\162BA.mp4"'
$ws.Range("B17").Value = 'Slime Squishing'
$ws.Range("C17").Value = '🤩 ASMR Slime Squishing so Satisfying Relaxing #SlimeSquishing #slimeasmr #asmrslime #satisfying'
$ws.Range("D17").Value = '🤩 ASMR Slime Squishing so Satisfying
#SlimeSquishing #Slime #ASMR #Satisfying #Shorts #YouTubeShorts
#pipingbags  #slimeasmr #asmrslime 
slime asmr, asmr slime, satisfying
Thank you for watching the video. Please subscribe to Slime Squishing to watch our latest videos! Wish you have relaxing moments with Slime Squishing.'
$ws.Range("E17").Value = '10:00'
$ws.Range("F17").Value = '''13/9/2025'
$ws.Range("G17").Value = 'Uploaded'

# Row 18
$ws.Range("A18").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Relaxing Squishy\130AA.mp4"'
$ws.Range("B18").Value = 'Relaxing SquishySquishy'
$ws.Range("C18").Value = '🍎🎈 Relaxing Squishy #bubblesound #stressballasmr #crunchystressballasmr #bubbleasmr'
$ws.Range("D18").Value = '#RelaxingSquishy #asmr #satisfying #relaxing #mixing #relax #stressrelief #bubblesound #stressballasmr #squishytungtungsahur #crunchystressballasmr #bubbleasmr
bubble sound, stress ball asmr, squishy tung tung sahur, crunchy stress ball asmr, bubble asmr
💬 What did you think of the video?
👇 Drop your thoughts in the comments below!
👍 If you enjoyed it, don’t forget to give it a like
🔔 Subscribe and turn on the bell so you never miss a new upload!
📢 Share this video with your friends too!
Thanks for watching and supporting Relaxing Squishy! 💖'
$ws.Range("E18").Value = '21:00'
$ws.Range("F18").Value = '''11/9/2025'

# Row 19
$ws.Range("A19").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Relaxing Squishy\132AA.mp4"'
$ws.Range("B19").Value = 'Relaxing SquishySquishy'
$ws.Range("C19").Value = '💥👉 Relaxing Squishy #bubblesound #stressballasmr #crunchystressballasmr #bubbleasmr
'
$ws.Range("D19").Value = '#RelaxingSquishy #asmr #satisfying #relaxing #mixing #relax #stressrelief #bubblesound #stressballasmr #squishytungtungsahur #crunchystressballasmr #bubbleasmr
bubble sound, stress ball asmr, squishy tung tung sahur, crunchy stress ball asmr, bubble asmr
💬 What did you think of the video?
👇 Drop your thoughts in the comments below!
👍 If you enjoyed it, don’t forget to give it a like
🔔 Subscribe and turn on the bell so you never miss a new upload!
📢 Share this video with your friends too!
Thanks for watching and supporting Relaxing Squishy! 💖'
$ws.Range("E19").Value = '1:00'
$ws.Range("F19").Value = '''12/9/2025'

# Row 20
$ws.Range("A20").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Relaxing Squishy\134AA.mp4"'
$ws.Range("B20").Value = 'Relaxing SquishySquishy'
$ws.Range("C20").Value = '🎉👉 Relaxing Squishy #bubblesound #stressballasmr #crunchystressballasmr #bubbleasmr
'
$ws.Range("D20").Value = '#RelaxingSquishy #asmr #satisfying #relaxing #mixing #relax #stressrelief #bubblesound #stressballasmr #squishytungtungsahur #crunchystressballasmr #bubbleasmr
bubble sound, stress ball asmr, squishy tung tung sahur, crunchy stress ball asmr, bubble asmr
💬 What did you think of the video?
👇 Drop your thoughts in the comments below!
👍 If you enjoyed it, don’t forget to give it a like
🔔 Subscribe and turn on the bell so you never miss a new upload!
📢 Share this video with your friends too!
Thanks for watching and supporting Relaxing Squishy! 💖'
$ws.Range("E20").Value = '9:00'
$ws.Range("F20").Value = '''12/9/2025'

# Row 21
$ws.Range("A21").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Relaxing Squishy\136AA.mp4"'
$ws.Range("B21").Value = 'Relaxing SquishySquishy'
$ws.Range("C21").Value = '💦🧁 Relaxing Squishy #bubblesound #stressballasmr #crunchystressballasmr #bubbleasmr
'
$ws.Range("D21").Value = '#RelaxingSquishy #asmr #satisfying #relaxing #mixing #relax #stressrelief #bubblesound #stressballasmr #squishytungtungsahur #crunchystressballasmr #bubbleasmr
bubble sound, stress ball asmr, squishy tung tung sahur, crunchy stress ball asmr, bubble asmr
💬 What did you think of the video?
👇 Drop your thoughts in the comments below!
👍 If you enjoyed it, don’t forget to give it a like
🔔 Subscribe and turn on the bell so you never miss a new upload!
📢 Share this video with your friends too!
Thanks for watching and supporting Relaxing Squishy! 💖'
$ws.Range("E21").Value = '21:00'
$ws.Range("F21").Value = '''12/9/2025'

# Row 22
$ws.Range("A22").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Relaxing Squishy\138AA.mp4"'
$ws.Range("B22").Value = 'Relaxing SquishySquishy'
$ws.Range("C22").Value = '👀 Relaxing Squishy #bubblesound #stressballasmr #crunchystressballasmr #bubbleasmr
'
$ws.Range("D22").Value = '#RelaxingSquishy #asmr #satisfying #relaxing #mixing #relax #stressrelief #bubblesound #stressballasmr #squishytungtungsahur #crunchystressballasmr #bubbleasmr
bubble sound, stress ball asmr, squishy tung tung sahur, crunchy stress ball asmr, bubble asmr
💬 What did you think of the video?
👇 Drop your thoughts in the comments below!
👍 If you enjoyed it, don’t forget to give it a like
🔔 Subscribe and turn on the bell so you never miss a new upload!
📢 Share this video with your friends too!
Thanks for watching and supporting Relaxing Squishy! 💖'
$ws.Range("E22").Value = '1:00'
$ws.Range("F22").Value = '''13/9/2025'

# Row 23
$ws.Range("A23").Value = '"\\192.168.1.92\Ổ Sever Mới\Định\Satisfy ASMR\New folder\45s - MUSIC\Tuấn\Relaxing Squishy\140AA.mp4"'
$ws.Range("B23").Value = 'Relaxing SquishySquishy'
$ws.Range("C23").Value = '💞💕 Relaxing Squishy #bubblesound #stressballasmr #crunchystressballasmr #bubbleasmr
'
$ws.Range("D23").Value = '#RelaxingSquishy #asmr #satisfying #relaxing #mixing #relax #stressrelief #bubblesound #stressballasmr #squishytungtungsahur #crunchystressballasmr #bubbleasmr
bubble sound, stress ball asmr, squishy tung tung sahur, crunchy stress ball asmr, bubble asmr
💬 What did you think of the video?
👇 Drop your thoughts in the comments below!
👍 If you enjoyed it, don’t forget to give it a like
🔔 Subscribe and turn on the bell so you never miss a new upload!
📢 Share this video with your friends too!
Thanks for watching and supporting Relaxing Squishy! 💖'
$ws.Range("E23").Value = '9:00'
$ws.Range("F23").Value = '''13/9/2025'
